$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scratch cell used to stage literal percent-like text ("NN%") so PasteSpecial(values-only)
# writes it without Excel's automatic text->percentage-number conversion or style churn.
$ws.Range("A1000").NumberFormat = "@"

$ws.Range("E2").Value = "2026-02-17 19:18:31"
$ws.Range("E3").Value = "2026-02-17 19:18:34"
$ws.Range("E4").Value = "2026-02-17 19:18:36"
$ws.Range("J4").Value = "1018.3 hPa"
$ws.Range("E5").Value = "2026-02-17 19:18:38"
$ws.Range("O5").Value = "-3.7 °C"
$ws.Range("E6").Value = "2026-02-17 19:18:41"
$ws.Range("J6").Value = "1018.3 hPa"
$ws.Range("E7").Value = "2026-02-17 19:18:43"
$ws.Range("J7").Value = "1018.2 hPa"
$ws.Range("E8").Value = "2026-02-17 19:18:46"
$ws.Range("J8").Value = "1018.1 hPa"
$ws.Range("E9").Value = "2026-02-17 19:18:48"
$ws.Range("A1000").Value = "56%"
$ws.Range("A1000").Copy()
$ws.Range("H9").PasteSpecial(-4163)
$ws.Range("O9").Value = "12.7 °C"
$ws.Range("E10").Value = "2026-02-17 19:18:50"
$ws.Range("A1000").Value = "73%"
$ws.Range("A1000").Copy()
$ws.Range("H10").PasteSpecial(-4163)
$ws.Range("K10").Value = "10.5 MJ/m2"
$ws.Range("O10").Value = "10.6 °C"
$ws.Range("E11").Value = "2026-02-17 19:18:53"
$ws.Range("A1000").Value = "48%"
$ws.Range("A1000").Copy()
$ws.Range("H11").PasteSpecial(-4163)
$ws.Range("O11").Value = "7.7 °C"
$ws.Range("E12").Value = "2026-02-17 19:18:55"
$ws.Range("A1000").Value = "58%"
$ws.Range("A1000").Copy()
$ws.Range("H12").PasteSpecial(-4163)
$ws.Range("O12").Value = "12.8 °C"
$ws.Range("E13").Value = "2026-02-17 19:18:57"
$ws.Range("A1000").Value = "43%"
$ws.Range("A1000").Copy()
$ws.Range("H13").PasteSpecial(-4163)
$ws.Range("E14").Value = "2026-02-17 19:19:00"
$ws.Range("A1000").Value = "66%"
$ws.Range("A1000").Copy()
$ws.Range("H14").PasteSpecial(-4163)
$ws.Range("E15").Value = "2026-02-17 19:19:02"
$ws.Range("A1000").Value = "56%"
$ws.Range("A1000").Copy()
$ws.Range("H15").PasteSpecial(-4163)
$ws.Range("O15").Value = "12.4 °C"
$ws.Range("E16").Value = "2026-02-17 19:19:04"
$ws.Range("E17").Value = "2026-02-17 19:19:07"
$ws.Range("A1000").Value = "80%"
$ws.Range("A1000").Copy()
$ws.Range("H17").PasteSpecial(-4163)
$ws.Range("E18").Value = "2026-02-17 19:19:09"
$ws.Range("J18").Value = "1018.5 hPa"
$ws.Range("E19").Value = "2026-02-17 19:19:12"
$ws.Range("A1000").Value = "73%"
$ws.Range("A1000").Copy()
$ws.Range("H19").PasteSpecial(-4163)
$ws.Range("E20").Value = "2026-02-17 19:19:14"
$ws.Range("A1000").Value = "62%"
$ws.Range("A1000").Copy()
$ws.Range("H20").PasteSpecial(-4163)
$ws.Range("E21").Value = "2026-02-17 19:19:17"
$ws.Range("A1000").Value = "37%"
$ws.Range("A1000").Copy()
$ws.Range("H21").PasteSpecial(-4163)
$ws.Range("J21").Value = "1016.9 hPa"
$ws.Range("E22").Value = "2026-02-17 19:19:19"
$ws.Range("E23").Value = "2026-02-17 19:19:21"
$ws.Range("A1000").Value = "70%"
$ws.Range("A1000").Copy()
$ws.Range("H23").PasteSpecial(-4163)
$ws.Range("I23").Value = "2.7 mm"
$ws.Range("E24").Value = "2026-02-17 19:19:23"
$ws.Range("O24").Value = "12.8 °C"
$ws.Range("E25").Value = "2026-02-17 19:19:26"
$ws.Range("A1000").Value = "51%"
$ws.Range("A1000").Copy()
$ws.Range("H25").PasteSpecial(-4163)
$ws.Range("E26").Value = "2026-02-17 19:19:28"
$ws.Range("E27").Value = "2026-02-17 19:19:30"
$ws.Range("A1000").Value = "52%"
$ws.Range("A1000").Copy()
$ws.Range("H27").PasteSpecial(-4163)
$ws.Range("E28").Value = "2026-02-17 19:19:33"
$ws.Range("J28").Value = "1018.1 hPa"
$ws.Range("E29").Value = "2026-02-17 19:19:35"
$ws.Range("A1000").Value = "64%"
$ws.Range("A1000").Copy()
$ws.Range("H29").PasteSpecial(-4163)
$ws.Range("N29").Value = "8.2 °C 18:49 TU"
$ws.Range("O29").Value = "12.3 °C"
$ws.Range("E30").Value = "2026-02-17 19:19:37"
$ws.Range("A1000").Value = "62%"
$ws.Range("A1000").Copy()
$ws.Range("H30").PasteSpecial(-4163)
$ws.Range("J30").Value = "1018.2 hPa"
$ws.Range("E31").Value = "2026-02-17 19:19:40"
$ws.Range("E32").Value = "2026-02-17 19:19:42"
$ws.Range("E33").Value = "2026-02-17 19:19:45"
$ws.Range("A1000").Value = "41%"
$ws.Range("A1000").Copy()
$ws.Range("H33").PasteSpecial(-4163)
$ws.Range("E34").Value = "2026-02-17 19:19:47"
$ws.Range("A1000").Value = "50%"
$ws.Range("A1000").Copy()
$ws.Range("H34").PasteSpecial(-4163)
$ws.Range("E35").Value = "2026-02-17 19:19:49"
$ws.Range("J35").Value = "1020.0 hPa"
$ws.Range("E36").Value = "2026-02-17 19:19:52"
$ws.Range("A1000").Value = "58%"
$ws.Range("A1000").Copy()
$ws.Range("H36").PasteSpecial(-4163)
$ws.Range("J36").Value = "1018.4 hPa"
$ws.Range("O36").Value = "12.5 °C"
$ws.Range("E37").Value = "2026-02-17 19:19:54"
$ws.Range("E38").Value = "2026-02-17 19:19:56"
$ws.Range("K38").Value = "10.7 MJ/m2"
$ws.Range("O38").Value = "11.2 °C"
$ws.Range("E39").Value = "2026-02-17 19:19:59"
$ws.Range("A1000").Value = "57%"
$ws.Range("A1000").Copy()
$ws.Range("H39").PasteSpecial(-4163)
$ws.Range("E40").Value = "2026-02-17 19:20:01"
$ws.Range("A1000").Value = "51%"
$ws.Range("A1000").Copy()
$ws.Range("H40").PasteSpecial(-4163)
$ws.Range("J40").Value = "1017.8 hPa"
$ws.Range("E41").Value = "2026-02-17 19:20:04"
$ws.Range("J41").Value = "1018.0 hPa"
$ws.Range("O41").Value = "16.8 °C"
$ws.Range("E42").Value = "2026-02-17 19:20:06"
$ws.Range("A1000").Value = "57%"
$ws.Range("A1000").Copy()
$ws.Range("H42").PasteSpecial(-4163)
$ws.Range("N42").Value = "9.8 °C 18:57 TU"
$ws.Range("E43").Value = "2026-02-17 19:20:08"
$ws.Range("E44").Value = "2026-02-17 19:20:11"
$ws.Range("I44").Value = "3.9 mm"
$ws.Range("O44").Value = "-3.1 °C"
$ws.Range("E45").Value = "2026-02-17 19:20:13"
$ws.Range("A1000").Value = "68%"
$ws.Range("A1000").Copy()
$ws.Range("H45").PasteSpecial(-4163)
$ws.Range("E46").Value = "2026-02-17 19:20:16"
$ws.Range("N46").Value = "12.5 °C 18:59 TU"
$ws.Range("O46").Value = "15.7 °C"

$ws.Range("A1000").Clear()
